# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed
# values, per commit "Updated cryptos list on Thu Jun  8 21:55:44 UTC 2023
# with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "26.603.82"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.853.44"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D4").Formula = "'1.002"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Formula = "'263.36"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D7").Formula = "'0.5255"
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("D8").Formula = "'0.3238"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").Formula = "'0.06782"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").Formula = "'18.92"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").Formula = "'0.7815"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Formula = "'0.07770"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "1.852.01"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Formula = "'13.98"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Formula = "'1.001"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Formula = "'0.000007958"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "26.642.79"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Formula = "'4.631"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").Formula = "'9.465"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Formula = "'6.008"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").Formula = "'143.59"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Formula = "'2.169"
$ws.Range("E25").Value = "  -6.79%  "
$ws.Range("D26").Formula = "'1.682"
$ws.Range("D27").Formula = "'17.04"
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").Formula = "'111.83"
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").Formula = "'4.182"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Formula = "'4.106"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Formula = "'1.131"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Formula = "'0.7195"
$ws.Range("E34").Value = "  +5.50%  "
$ws.Range("D35").Formula = "'2.869"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Formula = "'3.112"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Formula = "'2.252"
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("D38").Formula = "'0.01791"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Formula = "'0.4869"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Formula = "'0.9011"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").Formula = "'111.22"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Formula = "'5.960"
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("D43").Formula = "'1.001"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Formula = "'7.661"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Formula = "'0.05893"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Formula = "'8.990"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").Formula = "'0.8884"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("D51").Formula = "'59.82"
$ws.Range("E51").Value = "  +1.10%  "
